$wb = $excel.ActiveWorkbook
Write-Output $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
